$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (48) with the next sale record, mirroring the existing
# rows above it (same style/format as the rest of the "Automitivos"/"USA"
# "Limpa vidros" records already in the sheet).
$newRow = 48

$ws.Cells.Item($newRow, 1).Value = 45026
$ws.Cells.Item($newRow, 2).Value = 394
$ws.Cells.Item($newRow, 3).Value = 5
$ws.Cells.Item($newRow, 4).Value = "Automitivos"
$ws.Cells.Item($newRow, 5).Value = 639
$ws.Cells.Item($newRow, 6).Value = "USA"
$ws.Cells.Item($newRow, 7).Value = 8
$ws.Cells.Item($newRow, 8).Value = "Limpa vidros"
$ws.Cells.Item($newRow, 9).Value = 250
$ws.Cells.Item($newRow, 10).Value = 290
$ws.Cells.Item($newRow, 11).Value = 300
$ws.Cells.Item($newRow, 12).Value = 87000
$ws.Cells.Item($newRow, 13).Value = 19392

# Column A keeps the same date number format used by the rest of column A.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat
